# Add data for 2025-07-01
# Updates the 2025 (column L) year-to-date totals -- and, for a handful of
# late-reported records, the 2024 (column K) totals -- across the
# "Citywide Totals" summary sheet, the "By Neighborhood" summary sheet, and
# every affected per-neighborhood sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7901
$ws.Range("L2").Value = 3266
$ws.Range("L3").Value = 3372
$ws.Range("L4").Value = 844
$ws.Range("L6").Value = 2985
$ws.Range("K7").Value = 27557
$ws.Range("L7").Value = 10657

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 121

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 201
$ws.Range("L3").Value = 222
$ws.Range("L4").Value = 46
$ws.Range("L6").Value = 193
$ws.Range("L7").Value = 684

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 81
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 250

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 139
$ws.Range("L3").Value = 149
$ws.Range("L6").Value = 170
$ws.Range("L7").Value = 491

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 146

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L3").Value = 115
$ws.Range("L6").Value = 117
$ws.Range("L7").Value = 383

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L3").Value = 61
$ws.Range("L7").Value = 207

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 47
$ws.Range("L7").Value = 180

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 89
$ws.Range("L6").Value = 85
$ws.Range("L7").Value = 357
$ws.Range("L8").Value = 684
$ws.Range("L9").Value = 68
$ws.Range("L18").Value = 81
$ws.Range("L19").Value = 299
$ws.Range("L20").Value = 272
$ws.Range("L23").Value = 115
$ws.Range("L24").Value = 26
$ws.Range("L29").Value = 573
$ws.Range("L33").Value = 491
$ws.Range("L34").Value = 69
$ws.Range("L36").Value = 142
$ws.Range("L37").Value = 383
$ws.Range("L42").Value = 343
$ws.Range("L48").Value = 150
$ws.Range("L51").Value = 129
$ws.Range("L52").Value = 212
$ws.Range("L53").Value = 121
$ws.Range("L59").Value = 16
$ws.Range("K63").Value = 161
$ws.Range("L64").Value = 71
$ws.Range("L65").Value = 207
$ws.Range("L66").Value = 27
$ws.Range("L67").Value = 390
$ws.Range("L76").Value = 145
$ws.Range("L78").Value = 134
$ws.Range("L79").Value = 273
$ws.Range("L83").Value = 250
$ws.Range("L85").Value = 541
$ws.Range("L86").Value = 81
$ws.Range("L89").Value = 146
$ws.Range("L90").Value = 99
$ws.Range("L92").Value = 32
$ws.Range("L94").Value = 126
$ws.Range("L95").Value = 146
$ws.Range("L99").Value = 180
$ws.Range("K101").Value = 27557
$ws.Range("L101").Value = 10657

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 147
$ws.Range("L7").Value = 390

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L6").Value = 150
$ws.Range("L7").Value = 573

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L6").Value = 63
$ws.Range("L7").Value = 150

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 299

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L2").Value = 29
$ws.Range("L7").Value = 145

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 103
$ws.Range("L3").Value = 108
$ws.Range("L7").Value = 343

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L6").Value = 41
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L2").Value = 33
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 115

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 94
$ws.Range("L3").Value = 97
$ws.Range("L7").Value = 273

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 88
$ws.Range("L7").Value = 272

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L3").Value = 40
$ws.Range("L4").Value = 11
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 142

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 114
$ws.Range("L4").Value = 25
$ws.Range("L6").Value = 100
$ws.Range("L7").Value = 357

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L2").Value = 19
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 126

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L3").Value = 28
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 89

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 146

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L3").Value = 16
$ws.Range("L4").Value = 43
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L3").Value = 39
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 156
$ws.Range("L3").Value = 221
$ws.Range("L6").Value = 115
$ws.Range("L7").Value = 541

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 75
$ws.Range("L3").Value = 64
$ws.Range("L7").Value = 212
